# Apply a 3-way rotation of fields between rows 95, 97, 99 on the
# "Schedule" and "Registrar Schedule" sheets (the CS 396 / CS 398 /
# HNRS 251 course-offering rows got their data shuffled around), and
# bump the ExportTime timestamp on the "Metadata" sheet.
#
# Numeric-looking values (course numbers, loads, room capacities, ...)
# are stored as TEXT in this workbook, so they are written back with a
# leading apostrophe to stop Excel from auto-coercing them to numbers.
# Genuinely-numeric columns (MeetingDuration / Duration) are written as
# plain numbers.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, $value)

    $current = $range.Value()
    # Skip no-op writes so an unchanged cell (two rotated-from rows
    # happened to already share the same value) keeps its original
    # on-disk encoding/style untouched.
    if ("$current" -eq "$value") {
        return
    }

    if ($value -eq $null) {
        $range.Value = ""
    } elseif ($value -match '^-?[0-9]+(\.[0-9]+)?$') {
        # Numeric-looking text: force text storage with a quote-prefix.
        $range.Value = "'" + $value
    } else {
        $range.Value = $value
    }
}

function Set-NumValue {
    param($range, $value)

    $current = $range.Value()
    if ("$current" -eq "$value") {
        return
    }

    $range.Value = $value
}

# ---------------------------------------------------------------
# Sheet "Schedule": rotating fields live in columns C,E,F,H,I,L,M,N,O,P,R
# All of those are TEXT columns except N (MeetingDuration), which is
# a genuine number.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Schedule")

$textCols1 = @("C","E","F","H","I","L","M","O","P","R")
$numCols1  = @("N")

$before95 = @{}
$before97 = @{}
$before99 = @{}
foreach ($col in ($textCols1 + $numCols1)) {
    $before95[$col] = $ws1.Range("$col" + "95").Value()
    $before97[$col] = $ws1.Range("$col" + "97").Value()
    $before99[$col] = $ws1.Range("$col" + "99").Value()
}

# rotation: 95 <- 99, 97 <- 95, 99 <- 97
foreach ($col in $textCols1) {
    Set-TextValue $ws1.Range("$col" + "95") $before99[$col]
    Set-TextValue $ws1.Range("$col" + "97") $before95[$col]
    Set-TextValue $ws1.Range("$col" + "99") $before97[$col]
}
foreach ($col in $numCols1) {
    Set-NumValue $ws1.Range("$col" + "95") $before99[$col]
    Set-NumValue $ws1.Range("$col" + "97") $before95[$col]
    Set-NumValue $ws1.Range("$col" + "99") $before97[$col]
}

# ---------------------------------------------------------------
# Sheet "Registrar Schedule": mirrored layout, columns A,B,C,F,G,H,I,K,L,M,N
# All TEXT except L (Duration), a genuine number.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Registrar Schedule")

$textCols2 = @("A","B","C","F","G","H","I","K","M","N")
$numCols2  = @("L")

$rbefore95 = @{}
$rbefore97 = @{}
$rbefore99 = @{}
foreach ($col in ($textCols2 + $numCols2)) {
    $rbefore95[$col] = $ws2.Range("$col" + "95").Value()
    $rbefore97[$col] = $ws2.Range("$col" + "97").Value()
    $rbefore99[$col] = $ws2.Range("$col" + "99").Value()
}

foreach ($col in $textCols2) {
    Set-TextValue $ws2.Range("$col" + "95") $rbefore99[$col]
    Set-TextValue $ws2.Range("$col" + "97") $rbefore95[$col]
    Set-TextValue $ws2.Range("$col" + "99") $rbefore97[$col]
}
foreach ($col in $numCols2) {
    Set-NumValue $ws2.Range("$col" + "95") $rbefore99[$col]
    Set-NumValue $ws2.Range("$col" + "97") $rbefore95[$col]
    Set-NumValue $ws2.Range("$col" + "99") $rbefore97[$col]
}

# ---------------------------------------------------------------
# Sheet "Metadata": bump ExportTime
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Metadata")
Set-TextValue $ws3.Range("A2") "2025-04-17 11:54:18"
